$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$newValues = @{
    3 = 2078
    4 = 274
    5 = 72
    6 = 6372
    7 = 258
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
